# Fruta / hortaliza, semanal
# Insert one new weekly data row (Apio, "Primera") at the top of the data
# block (row 423), pushing all existing data rows down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 423; this pushes the former row 423
# (and everything below it) down to row 424, automatically growing the
# used range / dimension from R518 to R519.
$ws.Rows.Item(423).Insert()

# The row that used to be 423 is now 424. Copy over the descriptive
# (non price/date) fields from it into the new row 423, since the new
# record shares the same market / category / quality metadata.
$copyCols = @(1, 2, 3, 5, 6, 7, 8, 9, 14, 15, 17, 18)
foreach ($col in $copyCols) {
    $ws.Cells.Item(423, $col).Value = $ws.Cells.Item(424, $col).Value()
}

# Set the new record's own date and price figures.
$ws.Cells.Item(423, 4).Value = 44889    # D423 - Fecha
$ws.Cells.Item(423, 10).Value = 800     # J423 - Volumen
$ws.Cells.Item(423, 11).Value = 10000   # K423 - Precio minimo
$ws.Cells.Item(423, 12).Value = 11000   # L423 - Precio maximo
$ws.Cells.Item(423, 13).Value = 10500   # M423 - Precio promedio ponderado
$ws.Cells.Item(423, 16).Value = 1750    # P423 - Precio $/Kg
